# Refresh the "cryptos" list snapshot (values/percentages) as produced by
# the scheduled GitHub Actions update job.
# Price cells (column D) are forced to Text format before assignment so
# that values such as "1.00" / "0.800" / "9.70" keep their exact printed
# form instead of being auto-parsed into numbers (which would drop
# trailing zeros / introduce floating point noise). The style is reset
# back to Normal afterwards so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.313.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.734.54'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.97%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.14'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +11.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.268'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0637'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0896'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.980.31'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.736.28'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.27'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.563'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '28.299.13'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '242.37'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0756'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.38%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.66'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.70'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.84%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.65'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.55'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.91%  '
$ws.Range('E28').Value = '  +0.95%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  +3.17%  '
$ws.Range('E31').Value = '  +2.64%  '
$ws.Range('E32').Value = '  +1.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.505.99'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.27'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E35').Value = '  -2.06%  '
$ws.Range('E36').Value = '  +2.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.604'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.40'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.46%  '
$ws.Range('E39').Value = '  +2.10%  '
$ws.Range('E40').Value = '  +1.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '70.58'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.01%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.69'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.30'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.883.52'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.800'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('E47').Value = '  +8.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '91.05'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('E49').Value = '  +5.70%  '
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.19'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.33%  '
